$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.9611273080660837
$ws.Range("B3").Value = 0.03887269193391642
$ws.Range("B4").Value = 0.852891156462585
$ws.Range("B5").Value = 0.9704081632653061
$ws.Range("B6").Value = 0.8537029040016562
